# Applies a cyclic rotation of species-record data among rows 3, 4 and 5:
#   new row 3 <- old row 5
#   new row 4 <- old row 3
#   new row 5 <- old row 4
# Only the columns that actually carry per-record data are rotated
# (A, B, D, E, F, G, H, P, Q, R, AI); all other columns are identical
# across these three rows already, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","P","Q","R","AI")

# Capture the "before" values for rows 3, 4, 5
$data = @{}
foreach ($r in 3..5) {
    $data[$r] = @{}
    foreach ($col in $cols) {
        $data[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# New row assignment mapping (target row -> source row)
$mapping = @{ 3 = 5; 4 = 3; 5 = 4 }

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $data[$sourceRow][$col]
    }
}
